$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (GitHub Actions scheduled update).
# Numeric-looking price strings are written with a leading apostrophe and then
# restyled to Normal so Excel keeps them as plain text (matching the original
# inlineStr cells) instead of silently coercing them into Number cells.

$ws.Range("D2").Value = "70.715.48"
$ws.Range("E2").Value = "  -3.15%  "
$ws.Range("D3").Value = "3.824.54"
$ws.Range("E3").Value = "  -4.22%  "
$ws.Range("D5").Value = "'593.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.97%  "
$ws.Range("D6").Value = "'180.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.06%  "
$ws.Range("D7").Value = "'0.661"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "'0.750"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'0.175"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "'55.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.59%  "
$ws.Range("D12").Value = "'0.0000313"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'11.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "4.442.16"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").Value = "3.844.16"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "'20.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "'13.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.93%  "
$ws.Range("D18").Value = "'1.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.65%  "
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").Value = "70.622.96"
$ws.Range("E20").Value = "  -3.06%  "
$ws.Range("D21").Value = "'427.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").Value = "'4.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").Value = "'92.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.90%  "
$ws.Range("D24").Value = "'3.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.31%  "
$ws.Range("D25").Value = "'13.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.58%  "
$ws.Range("D26").Value = "'11.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").Value = "'3.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.01%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'10.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").Value = "'8.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +10.97%  "
$ws.Range("D31").Value = "'34.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").Value = "'13.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'46.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.06%  "
$ws.Range("D34").Value = "'0.124"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'635.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'67.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0964"
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("D38").Value = "'0.420"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.66%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.143"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").Value = "'3.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.02%  "
$ws.Range("D44").Value = "'0.0462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.55%  "
$ws.Range("D45").Value = "'9.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.45%  "
$ws.Range("D46").Value = "'2.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("D47").Value = "'0.141"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.65%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.33%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -15.75%  "
$ws.Range("D50").Value = "2.852.28"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "'0.000272"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
